$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 5714.316
$ws.Range("I28").Value = 5920.6665
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 5920.6665
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = -5435.6665
$ws.Range("N28").Value = -2970
$ws.Range("H51").Value = 4997.25
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516
$ws.Range("H61").Value = 792.3333
$ws.Range("I61").Value = 792.3333
$ws.Range("K61").Value = 2376.9999
$ws.Range("M61").Value = -2204.9999
$ws.Range("H70").Value = 3792.5625
$ws.Range("I70").Value = 4499.5
$ws.Range("J70").Value = 3691.5715
$ws.Range("K70").Value = 13498.5
$ws.Range("L70").Value = 11074.7145
$ws.Range("M70").Value = -13228.5
$ws.Range("N70").Value = -11614.7145
$ws.Range("H73").Value = 3792.5625
$ws.Range("I73").Value = 4499.5
$ws.Range("J73").Value = 3691.5715
$ws.Range("K73").Value = 13498.5
$ws.Range("L73").Value = 11074.7145
$ws.Range("M73").Value = -12562.5
$ws.Range("N73").Value = -12946.7145
$ws.Range("H80").Value = 1250.5
$ws.Range("I80").Value = 750
$ws.Range("J80").Value = 1751
$ws.Range("K80").Value = 2250
$ws.Range("L80").Value = 5253
$ws.Range("M80").Value = -1252
$ws.Range("N80").Value = -7249
$ws.Range("H83").Value = 1250.5
$ws.Range("I83").Value = 750
$ws.Range("J83").Value = 1751
$ws.Range("K83").Value = 6750
$ws.Range("L83").Value = 15759
$ws.Range("M83").Value = -1758
$ws.Range("N83").Value = -25743
$ws.Range("H86").Value = 12364
$ws.Range("I86").Value = 2193
$ws.Range("J86").Value = 20500.8
$ws.Range("K86").Value = 2193
$ws.Range("L86").Value = 20500.8
$ws.Range("M86").Value = -1070
$ws.Range("N86").Value = -22746.8
$ws.Range("H89").Value = 12364
$ws.Range("I89").Value = 2193
$ws.Range("J89").Value = 20500.8
$ws.Range("K89").Value = 10965
$ws.Range("L89").Value = 102504
$ws.Range("M89").Value = -5349
$ws.Range("N89").Value = -113736
$ws.Range("H92").Value = 986.8
$ws.Range("I92").Value = 944.6667
$ws.Range("K92").Value = 944.6667
$ws.Range("M92").Value = 303.3333
$ws.Range("H129").Value = 3690.6
$ws.Range("I129").Value = 1130.3334
$ws.Range("K129").Value = 3391.0002
$ws.Range("M129").Value = 1608.9998
$ws.Range("H138").Value = 3527.2727
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3527.2727
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10581.8181
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -20861.8181
$ws.Range("H141").Value = 5076.1113
$ws.Range("I141").Value = 5076.1113
$ws.Range("K141").Value = 15228.3339
$ws.Range("M141").Value = -10048.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 801.5
$ws.Range("I5").Value = 902
$ws.Range("K5").Value = 902
$ws.Range("M5").Value = -790
$ws.Range("H45").Value = 2913
$ws.Range("I45").Value = 2718.1667
$ws.Range("K45").Value = 2718.1667
$ws.Range("M45").Value = -2341.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 801.5
$ws.Range("I4").Value = 902
$ws.Range("K4").Value = 902
$ws.Range("M4").Value = -787
$ws.Range("H20").Value = 1849.5
$ws.Range("I20").Value = 1549.8334
$ws.Range("J20").Value = 2748.5
$ws.Range("K20").Value = 1549.8334
$ws.Range("L20").Value = 2748.5
$ws.Range("M20").Value = -1302.8334
$ws.Range("N20").Value = -3242.5
$ws.Range("H64").Value = 984.1111
$ws.Range("I64").Value = 1266.625
$ws.Range("J64").Value = 758.1
$ws.Range("K64").Value = 1266.625
$ws.Range("L64").Value = 758.1
$ws.Range("M64").Value = -1041.625
$ws.Range("N64").Value = -1208.1
$ws.Range("H67").Value = 984.1111
$ws.Range("I67").Value = 1266.625
$ws.Range("J67").Value = 758.1
$ws.Range("K67").Value = 1266.625
$ws.Range("L67").Value = 758.1
$ws.Range("M67").Value = -486.625
$ws.Range("N67").Value = -2318.1
$ws.Range("H86").Value = 2260.6875
$ws.Range("I86").Value = 2389.5833
$ws.Range("K86").Value = 2389.5833
$ws.Range("M86").Value = -1266.5833
$ws.Range("H89").Value = 2260.6875
$ws.Range("I89").Value = 2389.5833
$ws.Range("K89").Value = 11947.9165
$ws.Range("M89").Value = -6331.916499999999
$ws.Range("H94").Value = 2506.2
$ws.Range("I94").Value = 2506.6428
$ws.Range("K94").Value = 2506.6428
$ws.Range("M94").Value = -2055.6428
$ws.Range("H107").Value = 766.93335
$ws.Range("I107").Value = 377.72726
$ws.Range("J107").Value = 1837.25
$ws.Range("K107").Value = 377.72726
$ws.Range("L107").Value = 1837.25
$ws.Range("M107").Value = 1542.27274
$ws.Range("N107").Value = -5677.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3620
$ws.Range("J58").Value = 4333.3335
$ws.Range("L58").Value = 4333.3335
$ws.Range("N58").Value = -4739.3335
$ws.Range("H122").Value = 3338.75
$ws.Range("I122").Value = 4879.5
$ws.Range("J122").Value = 1798
$ws.Range("K122").Value = 14638.5
$ws.Range("L122").Value = 5394
$ws.Range("M122").Value = -12188.5
$ws.Range("N122").Value = -10294
$ws.Range("H136").Value = 3620
$ws.Range("J136").Value = 4333.3335
$ws.Range("L136").Value = 13000.0005
$ws.Range("N136").Value = -18100.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 1910.6
$ws.Range("J123").Value = 3277
$ws.Range("L123").Value = 9831
$ws.Range("N123").Value = -14731

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 43178.5
$ws.Range("J101").Value = 43178.5
$ws.Range("L101").Value = 43178.5
$ws.Range("N101").Value = -49668.5
$ws.Range("H107").Value = 300
$ws.Range("I107").Value = 300
$ws.Range("K107").Value = 300
$ws.Range("M107").Value = 1620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 26472.2
$ws.Range("J101").Value = 26472.2
$ws.Range("L101").Value = 26472.2
$ws.Range("N101").Value = -32962.2
$ws.Range("H122").Value = 3098.9092
$ws.Range("I122").Value = 2958.8
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 8876.400000000001
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -6426.400000000001
$ws.Range("N122").Value = -18400.0
$ws.Range("H132").Value = 21996.863
$ws.Range("J132").Value = 21748.75
$ws.Range("L132").Value = 65246.25
$ws.Range("N132").Value = -70306.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 17000
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 30000
$ws.Range("M32").Value = -3683
$ws.Range("N32").Value = -30634
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 353.5
$ws.Range("I107").Value = 359.44446
$ws.Range("K107").Value = 1078.33338
$ws.Range("M107").Value = 841.66662
$ws.Range("H122").Value = 1241.5
$ws.Range("I122").Value = 1275.3636
$ws.Range("J122").Value = 869
$ws.Range("K122").Value = 3826.0908
$ws.Range("L122").Value = 2607
$ws.Range("M122").Value = -1376.0908
$ws.Range("N122").Value = -7507.0
$ws.Range("H135").Value = 43799.8
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
